$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The dates in column A are stored as plain text. Excel's smart entry would
# reinterpret ambiguous "dd-mm-yyyy" strings (day <= 12) as real dates, so we
# force the cell to Text format while writing, then restore the default
# ("Normal") style so no stray formatting is left behind.
function Set-DateText($addr, $text) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 3: date separator change AND attendance counts updated (Invalid=1, Total=1)
Set-DateText "A3" "28-07-2022"
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

# Remaining rows: only the date separator changes from "/" to "-"
Set-DateText "A4"  "01-08-2022"
Set-DateText "A5"  "04-08-2022"
Set-DateText "A6"  "08-08-2022"
Set-DateText "A7"  "11-08-2022"
Set-DateText "A8"  "15-08-2022"
Set-DateText "A9"  "18-08-2022"
Set-DateText "A10" "22-08-2022"
Set-DateText "A11" "25-08-2022"
Set-DateText "A12" "29-08-2022"
Set-DateText "A13" "01-09-2022"
Set-DateText "A14" "05-09-2022"
Set-DateText "A15" "08-09-2022"
Set-DateText "A16" "12-09-2022"
Set-DateText "A17" "15-09-2022"
Set-DateText "A18" "19-09-2022"
Set-DateText "A19" "22-09-2022"
Set-DateText "A20" "26-09-2022"
Set-DateText "A21" "29-09-2022"
